$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 20:27"

# Update country stats rows (refreshed totals; some rows also involve a country-order change)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 7799116
$ws.Cells.Item(4, 3).Value = 22892
$ws.Cells.Item(4, 4).Value = 5000736
$ws.Cells.Item(4, 5).Value = 2581172
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 427
$ws.Cells.Item(4, 8).Value = 217208

$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 6902775
$ws.Cells.Item(5, 3).Value = 69787
$ws.Cells.Item(5, 4).Value = 5901710
$ws.Cells.Item(5, 5).Value = 894531
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 980
$ws.Cells.Item(5, 8).Value = 106534

$ws.Cells.Item(8, 1).Value = "España"
$ws.Cells.Item(8, 2).Value = 884381
$ws.Cells.Item(8, 3).Value = 5585
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 126
$ws.Cells.Item(8, 8).Value = 32688

$ws.Cells.Item(9, 1).Value = "Colombia"
$ws.Cells.Item(9, 2).Value = 877683
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 773973
$ws.Cells.Item(9, 5).Value = 76530
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 27180

$ws.Cells.Item(14, 1).Value = "Francia"
$ws.Cells.Item(14, 2).Value = 671638
$ws.Cells.Item(14, 3).Value = 18129
$ws.Cells.Item(14, 4).Value = 100306
$ws.Cells.Item(14, 5).Value = 538811
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 76
$ws.Cells.Item(14, 8).Value = 32521

$ws.Cells.Item(26, 1).Value = "Alemania"
$ws.Cells.Item(26, 2).Value = 313407
$ws.Cells.Item(26, 3).Value = 2294
$ws.Cells.Item(26, 4).Value = 267700
$ws.Cells.Item(26, 5).Value = 36042
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 13
$ws.Cells.Item(26, 8).Value = 9665

$ws.Cells.Item(32, 1).Value = "Ecuador"
$ws.Cells.Item(32, 2).Value = 145045
$ws.Cells.Item(32, 3).Value = 1514
$ws.Cells.Item(32, 4).Value = 120511
$ws.Cells.Item(32, 5).Value = 12393
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 398
$ws.Cells.Item(32, 8).Value = 12141

$ws.Cells.Item(33, 1).Value = "Marruecos"
$ws.Cells.Item(33, 2).Value = 142953
$ws.Cells.Item(33, 3).Value = 2929
$ws.Cells.Item(33, 4).Value = 120275
$ws.Cells.Item(33, 5).Value = 20192
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 47
$ws.Cells.Item(33, 8).Value = 2486

$ws.Cells.Item(54, 1).Value = "Etiopia"
$ws.Cells.Item(54, 2).Value = 81797
$ws.Cells.Item(54, 3).Value = 902
$ws.Cells.Item(54, 4).Value = 36434
$ws.Cells.Item(54, 5).Value = 44101
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 7
$ws.Cells.Item(54, 8).Value = 1262

$ws.Cells.Item(55, 1).Value = "Honduras"
$ws.Cells.Item(55, 2).Value = 81016
$ws.Cells.Item(55, 3).Value = 354
$ws.Cells.Item(55, 4).Value = 30590
$ws.Cells.Item(55, 5).Value = 47960
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 19
$ws.Cells.Item(55, 8).Value = 2466

$ws.Cells.Item(57, 1).Value = "Barein"
$ws.Cells.Item(57, 2).Value = 73932
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = 69411
$ws.Cells.Item(57, 5).Value = 4257
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 264

$ws.Cells.Item(58, 1).Value = "Uzbekistan"
$ws.Cells.Item(58, 2).Value = 60026
$ws.Cells.Item(58, 3).Value = 447
$ws.Cells.Item(58, 4).Value = 56837
$ws.Cells.Item(58, 5).Value = 2693
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 5
$ws.Cells.Item(58, 8).Value = 496

$ws.Cells.Item(59, 1).Value = "Moldavia"
$ws.Cells.Item(59, 2).Value = 59915
$ws.Cells.Item(59, 3).Value = 1121
$ws.Cells.Item(59, 4).Value = 43008
$ws.Cells.Item(59, 5).Value = 15483
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 18
$ws.Cells.Item(59, 8).Value = 1424

$ws.Cells.Item(64, 1).Value = "Argelia"
$ws.Cells.Item(64, 2).Value = 52658
$ws.Cells.Item(64, 3).Value = 138
$ws.Cells.Item(64, 4).Value = 36958
$ws.Cells.Item(64, 5).Value = 13917
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 12
$ws.Cells.Item(64, 8).Value = 1783

$ws.Cells.Item(68, 1).Value = "Ghana"
$ws.Cells.Item(68, 2).Value = 46947
$ws.Cells.Item(68, 3).Value = 118
$ws.Cells.Item(68, 4).Value = 46259
$ws.Cells.Item(68, 5).Value = 382
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 306

$ws.Cells.Item(74, 1).Value = "Irlanda"
$ws.Cells.Item(74, 2).Value = 40086
$ws.Cells.Item(74, 3).Value = 502
$ws.Cells.Item(74, 4).Value = 23364
$ws.Cells.Item(74, 5).Value = 14905
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 1817

$ws.Cells.Item(75, 1).Value = "Afganistan"
$ws.Cells.Item(75, 2).Value = 39616
$ws.Cells.Item(75, 3).Value = 68
$ws.Cells.Item(75, 4).Value = 33058
$ws.Cells.Item(75, 5).Value = 5088
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 1470

$ws.Cells.Item(156, 1).Value = "Sierra Leona"
$ws.Cells.Item(156, 2).Value = 2293
$ws.Cells.Item(156, 3).Value = 6
$ws.Cells.Item(156, 4).Value = 1721
$ws.Cells.Item(156, 5).Value = 500
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 72

$ws.Cells.Item(175, 1).Value = "Curazao"
$ws.Cells.Item(175, 2).Value = 532
$ws.Cells.Item(175, 3).Value = 27
$ws.Cells.Item(175, 4).Value = 281
$ws.Cells.Item(175, 5).Value = 250
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 1

$ws.Cells.Item(176, 1).Value = "Taiwan"
$ws.Cells.Item(176, 2).Value = 524
$ws.Cells.Item(176, 3).Value = 1
$ws.Cells.Item(176, 4).Value = 485
$ws.Cells.Item(176, 5).Value = 32
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 7

$ws.Cells.Item(177, 1).Value = "Burundi"
$ws.Cells.Item(177, 2).Value = 515
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 472
$ws.Cells.Item(177, 5).Value = 42
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 1

$ws.Cells.Item(178, 1).Value = "Tanzania"
$ws.Cells.Item(178, 2).Value = 509
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 183
$ws.Cells.Item(178, 5).Value = 305
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 21

$ws.Cells.Item(189, 1).Value = "Monaco"
$ws.Cells.Item(189, 2).Value = 229
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 207
$ws.Cells.Item(189, 5).Value = 20
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 2

